# Update outputs-r202 (previous copy of ful-path.csv):
# the quadratic-svm-score values for every genome in column B
# (rows 2-32) change from 1 -> 0. Column A keeps its existing
# text-formatted style (re-applied here to mirror the style churn
# recorded in the source workbook), column C ("prediction") is
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Re-assert the text number format on every cell that already carried
# it: the whole header row (A1:C1) plus the "Row" column A2:A32.
# Style churn in the source file re-allocates these onto a new (but
# functionally identical) style slot.
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A2:A" + $lastRow).NumberFormat = "@"

# The actual data edit: every score in column B (rows 2..lastRow)
# moves from 1 to 0.
$ws.Range("B2:B" + $lastRow).Value = 0
